$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -13.376
$ws.Range("B9").Value = 5.312
$ws.Range("C9").Value = -10.623
$ws.Range("C11").Value = -12.818
$ws.Range("B18").Value = 5.275
$ws.Range("B20").Value = 6.601999999999999
$ws.Range("C23").Value = -13.075
$ws.Range("C24").Value = -12.849
$ws.Range("C26").Value = -12.807
$ws.Range("B27").Value = 5.715
$ws.Range("C34").Value = -12.148
$ws.Range("B35").Value = 8.313000000000001
$ws.Range("C35").Value = -12.373
$ws.Range("C48").Value = -11.861
$ws.Range("C49").Value = -13.285
$ws.Range("C52").Value = -11.743
$ws.Range("C66").Value = -11.574
$ws.Range("C67").Value = -10.886
$ws.Range("B69").Value = 5.827
$ws.Range("B76").Value = 6.11
$ws.Range("B78").Value = 8.609999999999999
$ws.Range("C78").Value = -11.81
$ws.Range("C80").Value = -12.21
$ws.Range("B82").Value = 5.486
$ws.Range("B83").Value = 5.598
$ws.Range("B93").Value = 4.973000000000001
$ws.Range("C99").Value = -12.048
$ws.Range("C104").Value = -12.907

$wb.Save()
